$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.418.50'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.805.27'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.86'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4517'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.47%  '

$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.20'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07059'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8902'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +2.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07820'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("E13").Value = '  +0.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.838.55'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.273'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.298'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.33'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -1.42%  '

$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008489'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.008'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.449.20'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -0.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.17'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.958'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.063.87'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +1.61%  '

$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.962'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.27'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.069'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +4.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.07'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -0.57%  '

$ws.Range("E31").Value = '  -0.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08686'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.077'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -1.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.794'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +11.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.470'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7264'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.107'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.079'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("E39").Value = '  +1.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.913'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +1.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05111'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5056'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +2.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.793'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("E44").Value = '  -3.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.019'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.009'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4667'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +1.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.979'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.08'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -1.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.574'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05970'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -0.42%  '
